$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Fix D993:D1010: convert BSE code cells from text to numeric ---
for ($r = 993; $r -le 1010; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = $cell.Value()
}

# --- Append new rows 1011:1042 (stock.yaml break-out) ---
$ws.Cells.Item(1011, 1).Value = 1
$ws.Cells.Item(1011, 2).Value = "ABB"
$ws.Cells.Item(1011, 3).Value = "Abb India Limited"
$ws.Cells.Item(1011, 4).Value = 500002
$ws.Cells.Item(1011, 5).Value = 0.62
$ws.Cells.Item(1011, 6).Value = 7537.55
$ws.Cells.Item(1011, 7).Value = 339612
$ws.Cells.Item(1011, 8).Value = "day"
$ws.Cells.Item(1011, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1012, 1).Value = 2
$ws.Cells.Item(1012, 2).Value = "POLYCAB"
$ws.Cells.Item(1012, 3).Value = "Polycab India Ltd"
$ws.Cells.Item(1012, 4).Value = 542652
$ws.Cells.Item(1012, 5).Value = 0.82
$ws.Cells.Item(1012, 6).Value = 7412.5
$ws.Cells.Item(1012, 7).Value = 270695
$ws.Cells.Item(1012, 8).Value = "day"
$ws.Cells.Item(1012, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1013, 1).Value = 3
$ws.Cells.Item(1013, 2).Value = "ATUL"
$ws.Cells.Item(1013, 3).Value = "Atul Limited"
$ws.Cells.Item(1013, 4).Value = 500027
$ws.Cells.Item(1013, 5).Value = 0.63
$ws.Cells.Item(1013, 6).Value = 7399.95
$ws.Cells.Item(1013, 7).Value = 42281
$ws.Cells.Item(1013, 8).Value = "day"
$ws.Cells.Item(1013, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1014, 1).Value = 4
$ws.Cells.Item(1014, 2).Value = "TATAELXSI"
$ws.Cells.Item(1014, 3).Value = "Tata Elxsi Limited"
$ws.Cells.Item(1014, 4).Value = 500408
$ws.Cells.Item(1014, 5).Value = 5.74
$ws.Cells.Item(1014, 6).Value = 7098.55
$ws.Cells.Item(1014, 7).Value = 1080035
$ws.Cells.Item(1014, 8).Value = "day"
$ws.Cells.Item(1014, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1015, 1).Value = 5
$ws.Cells.Item(1015, 2).Value = "MCX"
$ws.Cells.Item(1015, 3).Value = "Multi Commodity Exchange Of India Limited"
$ws.Cells.Item(1015, 4).Value = 534091
$ws.Cells.Item(1015, 5).Value = -1.68
$ws.Cells.Item(1015, 6).Value = 6216.5
$ws.Cells.Item(1015, 7).Value = 262434
$ws.Cells.Item(1015, 8).Value = "day"
$ws.Cells.Item(1015, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1016, 1).Value = 6
$ws.Cells.Item(1016, 2).Value = "JKCEMENT"
$ws.Cells.Item(1016, 3).Value = "Jk Cement Limited"
$ws.Cells.Item(1016, 4).Value = 532644
$ws.Cells.Item(1016, 5).Value = 3.38
$ws.Cells.Item(1016, 6).Value = 4647.4
$ws.Cells.Item(1016, 7).Value = 264154
$ws.Cells.Item(1016, 8).Value = "day"
$ws.Cells.Item(1016, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1017, 1).Value = 7
$ws.Cells.Item(1017, 2).Value = "KEI"
$ws.Cells.Item(1017, 3).Value = "Kei Industries Limited"
$ws.Cells.Item(1017, 4).Value = 517569
$ws.Cells.Item(1017, 5).Value = -0.34
$ws.Cells.Item(1017, 6).Value = 4340.25
$ws.Cells.Item(1017, 7).Value = 274021
$ws.Cells.Item(1017, 8).Value = "day"
$ws.Cells.Item(1017, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1018, 1).Value = 8
$ws.Cells.Item(1018, 2).Value = "LALPATHLAB"
$ws.Cells.Item(1018, 3).Value = "Dr. Lal Path Labs Ltd."
$ws.Cells.Item(1018, 4).Value = 539524
$ws.Cells.Item(1018, 5).Value = -1.5
$ws.Cells.Item(1018, 6).Value = 2972.9
$ws.Cells.Item(1018, 7).Value = 139632
$ws.Cells.Item(1018, 8).Value = "day"
$ws.Cells.Item(1018, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1019, 1).Value = 9
$ws.Cells.Item(1019, 2).Value = "HINDUNILVR"
$ws.Cells.Item(1019, 3).Value = "Hindustan Unilever Limited"
$ws.Cells.Item(1019, 4).Value = 500696
$ws.Cells.Item(1019, 5).Value = 0.15
$ws.Cells.Item(1019, 6).Value = 2482.85
$ws.Cells.Item(1019, 7).Value = 1641588
$ws.Cells.Item(1019, 8).Value = "day"
$ws.Cells.Item(1019, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1020, 1).Value = 10
$ws.Cells.Item(1020, 2).Value = "INDIAMART"
$ws.Cells.Item(1020, 3).Value = "Indiamart Intermesh Ltd"
$ws.Cells.Item(1020, 4).Value = 542726
$ws.Cells.Item(1020, 5).Value = -0.22
$ws.Cells.Item(1020, 6).Value = 2357.55
$ws.Cells.Item(1020, 7).Value = 143189
$ws.Cells.Item(1020, 8).Value = "day"
$ws.Cells.Item(1020, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1021, 1).Value = 11
$ws.Cells.Item(1021, 2).Value = "NESTLEIND"
$ws.Cells.Item(1021, 3).Value = "Nestle India Limited"
$ws.Cells.Item(1021, 4).Value = 500790
$ws.Cells.Item(1021, 5).Value = 0.44
$ws.Cells.Item(1021, 6).Value = 2261.7
$ws.Cells.Item(1021, 7).Value = 766232
$ws.Cells.Item(1021, 8).Value = "day"
$ws.Cells.Item(1021, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1022, 1).Value = 12
$ws.Cells.Item(1022, 2).Value = "CYIENT"
$ws.Cells.Item(1022, 3).Value = "Cyient Limited"
$ws.Cells.Item(1022, 4).Value = 532175
$ws.Cells.Item(1022, 5).Value = 0.35
$ws.Cells.Item(1022, 6).Value = 1876.6
$ws.Cells.Item(1022, 7).Value = 292687
$ws.Cells.Item(1022, 8).Value = "day"
$ws.Cells.Item(1022, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1023, 1).Value = 13
$ws.Cells.Item(1023, 2).Value = "ASTRAL"
$ws.Cells.Item(1023, 3).Value = "Astral Poly Technik Limited"
$ws.Cells.Item(1023, 4).Value = 532830
$ws.Cells.Item(1023, 5).Value = 1.62
$ws.Cells.Item(1023, 6).Value = 1845.7
$ws.Cells.Item(1023, 7).Value = 624866
$ws.Cells.Item(1023, 8).Value = "day"
$ws.Cells.Item(1023, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1024, 1).Value = 14
$ws.Cells.Item(1024, 2).Value = "PVRINOX"
$ws.Cells.Item(1024, 3).Value = "PVR Inox Ltd"
$ws.Cells.Item(1024, 4).Value = 532689
$ws.Cells.Item(1024, 5).Value = 1.64
$ws.Cells.Item(1024, 6).Value = 1598.3
$ws.Cells.Item(1024, 7).Value = 583115
$ws.Cells.Item(1024, 8).Value = "day"
$ws.Cells.Item(1024, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1025, 1).Value = 15
$ws.Cells.Item(1025, 2).Value = "BHARATFORG"
$ws.Cells.Item(1025, 3).Value = "Bharat Forge Limited"
$ws.Cells.Item(1025, 4).Value = 500493
$ws.Cells.Item(1025, 5).Value = 0.49
$ws.Cells.Item(1025, 6).Value = 1350.4
$ws.Cells.Item(1025, 7).Value = 791839
$ws.Cells.Item(1025, 8).Value = "day"
$ws.Cells.Item(1025, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1026, 1).Value = 16
$ws.Cells.Item(1026, 2).Value = "MGL"
$ws.Cells.Item(1026, 3).Value = "Mahanagar Gas Limited"
$ws.Cells.Item(1026, 4).Value = 539957
$ws.Cells.Item(1026, 5).Value = 3.38
$ws.Cells.Item(1026, 6).Value = 1278.25
$ws.Cells.Item(1026, 7).Value = 1609969
$ws.Cells.Item(1026, 8).Value = "day"
$ws.Cells.Item(1026, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1027, 1).Value = 17
$ws.Cells.Item(1027, 2).Value = "DRREDDY"
$ws.Cells.Item(1027, 3).Value = "Dr. Reddy's Laboratories Limited"
$ws.Cells.Item(1027, 4).Value = 500124
$ws.Cells.Item(1027, 5).Value = 0.23
$ws.Cells.Item(1027, 6).Value = 1224.5
$ws.Cells.Item(1027, 7).Value = 2237453
$ws.Cells.Item(1027, 8).Value = "day"
$ws.Cells.Item(1027, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1028, 1).Value = 18
$ws.Cells.Item(1028, 2).Value = "SYNGENE"
$ws.Cells.Item(1028, 3).Value = "Syngene International Limited"
$ws.Cells.Item(1028, 4).Value = 539268
$ws.Cells.Item(1028, 5).Value = -1.37
$ws.Cells.Item(1028, 6).Value = 934.85
$ws.Cells.Item(1028, 7).Value = 623487
$ws.Cells.Item(1028, 8).Value = "day"
$ws.Cells.Item(1028, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1029, 1).Value = 19
$ws.Cells.Item(1029, 2).Value = "TATAMOTORS"
$ws.Cells.Item(1029, 3).Value = "Tata Motors Limited"
$ws.Cells.Item(1029, 4).Value = 500570
$ws.Cells.Item(1029, 5).Value = 1.42
$ws.Cells.Item(1029, 6).Value = 801.25
$ws.Cells.Item(1029, 7).Value = 9874135
$ws.Cells.Item(1029, 8).Value = "day"
$ws.Cells.Item(1029, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1030, 1).Value = 20
$ws.Cells.Item(1030, 2).Value = "MARICO"
$ws.Cells.Item(1030, 3).Value = "Marico Limited"
$ws.Cells.Item(1030, 4).Value = 531642
$ws.Cells.Item(1030, 5).Value = -0.75
$ws.Cells.Item(1030, 6).Value = 641.7
$ws.Cells.Item(1030, 7).Value = 1018776
$ws.Cells.Item(1030, 8).Value = "day"
$ws.Cells.Item(1030, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1031, 1).Value = 21
$ws.Cells.Item(1031, 2).Value = "GRANULES"
$ws.Cells.Item(1031, 3).Value = "Granules India Limited"
$ws.Cells.Item(1031, 4).Value = 532482
$ws.Cells.Item(1031, 5).Value = -10.23
$ws.Cells.Item(1031, 6).Value = 534.15
$ws.Cells.Item(1031, 7).Value = 19646137
$ws.Cells.Item(1031, 8).Value = "day"
$ws.Cells.Item(1031, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1032, 1).Value = 22
$ws.Cells.Item(1032, 2).Value = "DABUR"
$ws.Cells.Item(1032, 3).Value = "Dabur India Limited"
$ws.Cells.Item(1032, 4).Value = 500096
$ws.Cells.Item(1032, 5).Value = -0.26
$ws.Cells.Item(1032, 6).Value = 522.55
$ws.Cells.Item(1032, 7).Value = 3691100
$ws.Cells.Item(1032, 8).Value = "day"
$ws.Cells.Item(1032, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1033, 1).Value = 23
$ws.Cells.Item(1033, 2).Value = "GUJGASLTD"
$ws.Cells.Item(1033, 3).Value = "Gujarat Gas Limited"
$ws.Cells.Item(1033, 4).Value = 539336
$ws.Cells.Item(1033, 5).Value = 3.6
$ws.Cells.Item(1033, 6).Value = 497.8
$ws.Cells.Item(1033, 7).Value = 1356678
$ws.Cells.Item(1033, 8).Value = "day"
$ws.Cells.Item(1033, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1034, 1).Value = 24
$ws.Cells.Item(1034, 2).Value = "DELHIVERY"
$ws.Cells.Item(1034, 3).Value = "Delhivery Ltd"
$ws.Cells.Item(1034, 4).Value = 543529
$ws.Cells.Item(1034, 5).Value = 2.36
$ws.Cells.Item(1034, 6).Value = 341.15
$ws.Cells.Item(1034, 7).Value = 4541844
$ws.Cells.Item(1034, 8).Value = "day"
$ws.Cells.Item(1034, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1035, 1).Value = 25
$ws.Cells.Item(1035, 2).Value = "PETRONET"
$ws.Cells.Item(1035, 3).Value = "Petronet Lng Limited"
$ws.Cells.Item(1035, 4).Value = 532522
$ws.Cells.Item(1035, 5).Value = 0.22
$ws.Cells.Item(1035, 6).Value = 339.3
$ws.Cells.Item(1035, 7).Value = 1793652
$ws.Cells.Item(1035, 8).Value = "day"
$ws.Cells.Item(1035, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1036, 1).Value = 26
$ws.Cells.Item(1036, 2).Value = "M&MFIN"
$ws.Cells.Item(1036, 3).Value = "Mahindra & Mahindra Financial Services Limited"
$ws.Cells.Item(1036, 4).Value = 532720
$ws.Cells.Item(1036, 5).Value = 2.85
$ws.Cells.Item(1036, 6).Value = 285.3
$ws.Cells.Item(1036, 7).Value = 3319636
$ws.Cells.Item(1036, 8).Value = "day"
$ws.Cells.Item(1036, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1037, 1).Value = 27
$ws.Cells.Item(1037, 2).Value = "HUDCO"
$ws.Cells.Item(1037, 3).Value = "Housing and Urban Development Corporation"
$ws.Cells.Item(1037, 4).Value = 540530
$ws.Cells.Item(1037, 5).Value = -1.03
$ws.Cells.Item(1037, 6).Value = 235.54
$ws.Cells.Item(1037, 7).Value = 5281185
$ws.Cells.Item(1037, 8).Value = "day"
$ws.Cells.Item(1037, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1038, 1).Value = 28
$ws.Cells.Item(1038, 2).Value = "FEDERALBNK"
$ws.Cells.Item(1038, 3).Value = "The Federal Bank  Limited"
$ws.Cells.Item(1038, 4).Value = 500469
$ws.Cells.Item(1038, 5).Value = 0.42
$ws.Cells.Item(1038, 6).Value = 209.96
$ws.Cells.Item(1038, 7).Value = 6357323
$ws.Cells.Item(1038, 8).Value = "day"
$ws.Cells.Item(1038, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1039, 1).Value = 29
$ws.Cells.Item(1039, 2).Value = "IRFC"
$ws.Cells.Item(1039, 3).Value = "Indian Railway Finance Corporation Ltd"
$ws.Cells.Item(1039, 4).Value = 543257
$ws.Cells.Item(1039, 5).Value = 0.62
$ws.Cells.Item(1039, 6).Value = 148.2
$ws.Cells.Item(1039, 7).Value = 12162775
$ws.Cells.Item(1039, 8).Value = "day"
$ws.Cells.Item(1039, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1040, 1).Value = 30
$ws.Cells.Item(1040, 2).Value = "UNIONBANK"
$ws.Cells.Item(1040, 3).Value = "Union Bank Of India"
$ws.Cells.Item(1040, 4).Value = 532477
$ws.Cells.Item(1040, 5).Value = 5.29
$ws.Cells.Item(1040, 6).Value = 126.06
$ws.Cells.Item(1040, 7).Value = 28256395
$ws.Cells.Item(1040, 8).Value = "day"
$ws.Cells.Item(1040, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1041, 1).Value = 31
$ws.Cells.Item(1041, 2).Value = "SJVN"
$ws.Cells.Item(1041, 3).Value = "Sjvn Limited"
$ws.Cells.Item(1041, 4).Value = 533206
$ws.Cells.Item(1041, 5).Value = 2.83
$ws.Cells.Item(1041, 6).Value = 117.03
$ws.Cells.Item(1041, 7).Value = 6264268
$ws.Cells.Item(1041, 8).Value = "day"
$ws.Cells.Item(1041, 9).Value = "03/12/2024 11:35:18"

$ws.Cells.Item(1042, 1).Value = 32
$ws.Cells.Item(1042, 2).Value = "NHPC"
$ws.Cells.Item(1042, 3).Value = "Nhpc Limited"
$ws.Cells.Item(1042, 4).Value = 533098
$ws.Cells.Item(1042, 5).Value = 0.61
$ws.Cells.Item(1042, 6).Value = 81.98
$ws.Cells.Item(1042, 7).Value = 13912781
$ws.Cells.Item(1042, 8).Value = "day"
$ws.Cells.Item(1042, 9).Value = "03/12/2024 11:35:18"
